$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 247, shifting rows 247:350 down to 248:351.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new record.
$ws.Cells.Item(247, 1).Value = 6
$ws.Cells.Item(247, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(247, 3).Value = "Metropolitana"
$ws.Cells.Item(247, 4).Value = 44900
$ws.Cells.Item(247, 5).Value = 13
$ws.Cells.Item(247, 6).Value = 100112026
$ws.Cells.Item(247, 7).Value = "Haba"
$ws.Cells.Item(247, 8).Value = "Sin especificar"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 400
$ws.Cells.Item(247, 11).Value = 8000
$ws.Cells.Item(247, 12).Value = 9000
$ws.Cells.Item(247, 13).Value = 8425
$ws.Cells.Item(247, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(247, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(247, 16).Value = 337
$ws.Cells.Item(247, 17).Value = 25
$ws.Cells.Item(247, 18).Value = "Hortaliza"
